$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix typo in row 67 title: "Mo4" -> "M04" ---
$ws.Range("C67").Value = "Banker's Algorithm Question | Operating System - M04 P06"

# --- Add new row to the log table (Post 60) ---
$tbl = $ws.ListObjects.Item("Table2")
$newRow = $tbl.ListRows.Add()

$ws.Range("B70").Value = 60
$ws.Range("C70").Value = "Question on Deadlock | Operating System - M04 P07"

$ws.Range("D70").Value = 44178
$ws.Range("D70").NumberFormat = "m/d/yy"

$ws.Range("E70").Value = "https://programmingport.hashnode.dev/question-on-deadlock-or-operating-system-m04-p07"
$ws.Range("E70").Style = "Hyperlink"

$ws.Range("F70").Value = "https://dev.to/rahulmishra05/question-on-deadlock-operating-system-m04-p07-37hh"
$ws.Range("F70").Style = "Hyperlink"

# --- Match the workbook's last-saved selection/viewport ---
$ws.Range("E70").Select() | Out-Null
